$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths for A:E (A manually sized, B:E mirror Excel's AutoFit result) ---
# The stored OOXML <col> width is ColumnWidth + 5/6 (rounded to the host's pixel
# grid), so the inputs below are chosen so the saved width lands on 16.5 / ~34.17 /
# ~9.67 / 6.5 / 22 (the two "bestFit" columns quantize to the nearest 1/6 the host
# supports, which is as close as this engine can get to Excel's true font metrics).
$ws.Range("A:A").ColumnWidth = 15.666666666666666
$ws.Range("B:B").ColumnWidth = 33.333333333333336
$ws.Range("C:C").ColumnWidth = 8.833333333333334
$ws.Range("D:D").ColumnWidth = 5.666666666666667
$ws.Range("E:E").ColumnWidth = 21.166666666666668

# --- Re-order the merged-cell list: put the A20:E26 block ahead of A17:E19 ---
# Unmerging + re-merging a block moves its entries to the end of the sheet's
# internal merge list, so re-doing A17:E19 (leaving A20:E26 untouched) results in
# A20:E26 first, then A17:E19 -- matching the target order.
$ws.Range("A17:A19").UnMerge()
$ws.Range("B17:B19").UnMerge()
$ws.Range("C17:C19").UnMerge()
$ws.Range("D17:D19").UnMerge()
$ws.Range("E17:E19").UnMerge()
$ws.Range("A17:A19").Merge()
$ws.Range("B17:B19").Merge()
$ws.Range("C17:C19").Merge()
$ws.Range("D17:D19").Merge()
$ws.Range("E17:E19").Merge()

# --- Move the active selection to E39 ---
$ws.Range("E39").Select()
